$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Amy and her family are going on a picnic. They pack sandwiches, fruits, and drinks in a basket. Amy helps her mom make the sandwiches. They also bring a blanket to sit on. When they arrive at the park, Amy and her brother play tag while their parents set up the picnic. They all enjoy the delicious food and the beautiful weather."
$ws.Range("C2").Value = "Amy e sua família vão fazer um piquenique. Eles preparam sanduíches, frutas e bebidas em uma cesta. Amy ajuda sua mãe a fazer os sanduíches. Eles também trazem um cobertor para sentar. Quando chegam ao parque, Amy e seu irmão brincam de pega enquanto seus pais arrumam o piquenique. Todos aproveitam a comida deliciosa e o clima agradável."
$ws.Range("E2").Value = 73
$ws.Range("F2").Value = "/OK"
